$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 268.33334
$ws.Range("I2").Value = 306.66666
$ws.Range("J2").Value = 230
$ws.Range("K2").Value = 306.66666
$ws.Range("L2").Value = 230
$ws.Range("M2").Value = -193.66666
$ws.Range("N2").Value = -456
$ws.Range("H21").Value = 23000
$ws.Range("J21").Value = 23000
$ws.Range("L21").Value = 23000
$ws.Range("N21").Value = -23936
$ws.Range("H23").Value = 23000
$ws.Range("J23").Value = 23000
$ws.Range("L23").Value = 23000
$ws.Range("N23").Value = -23468
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 1500
$ws.Range("M29").Value = -1219
$ws.Range("H32").Value = 1250.2778
$ws.Range("I32").Value = 1300.3334
$ws.Range("J32").Value = 1240.2667
$ws.Range("K32").Value = 1300.3334
$ws.Range("L32").Value = 1240.2667
$ws.Range("M32").Value = -974.3334
$ws.Range("N32").Value = -1892.2667
$ws.Range("H38").Value = 450.5
$ws.Range("I38").Value = 378.33334
$ws.Range("J38").Value = 1100
$ws.Range("K38").Value = 1135.00002
$ws.Range("L38").Value = 3300
$ws.Range("M38").Value = -763.0000199999999
$ws.Range("N38").Value = -4044
$ws.Range("H58").Value = 408.75
$ws.Range("I58").Value = 154
$ws.Range("J58").Value = 833.3333
$ws.Range("K58").Value = 462
$ws.Range("L58").Value = 2499.9999
$ws.Range("M58").Value = -312
$ws.Range("N58").Value = -2799.9999
$ws.Range("H87").Value = 57687.332
$ws.Range("J87").Value = 57687.332
$ws.Range("L87").Value = 57687.332
$ws.Range("N87").Value = -60183.332
$ws.Range("H90").Value = 57687.332
$ws.Range("J90").Value = 57687.332
$ws.Range("L90").Value = 173061.996
$ws.Range("N90").Value = -185541.996
$ws.Range("H93").Value = 35601
$ws.Range("J93").Value = 35601
$ws.Range("L93").Value = 35601
$ws.Range("N93").Value = -40593
$ws.Range("H98").Value = 2148.1765
$ws.Range("J98").Value = 1029.6666
$ws.Range("L98").Value = 1029.6666
$ws.Range("N98").Value = -4025.6666
$ws.Range("H122").Value = 2148.1765
$ws.Range("J122").Value = 1029.6666
$ws.Range("L122").Value = 3088.9998
$ws.Range("N122").Value = -7988.9998
$ws.Range("H138").Value = 1932.9131
$ws.Range("I138").Value = 1130.5
$ws.Range("J138").Value = 3181.111
$ws.Range("K138").Value = 3391.5
$ws.Range("L138").Value = 9543.332999999999
$ws.Range("M138").Value = 1748.5
$ws.Range("N138").Value = -19823.333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1697.8462
$ws.Range("I2").Value = 1108
$ws.Range("K2").Value = 1108
$ws.Range("M2").Value = -995
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 12
$ws.Range("K6").Value = 12
$ws.Range("M6").Value = 161
$ws.Range("H116").Value = 1697.8462
$ws.Range("I116").Value = 1108
$ws.Range("K116").Value = 1108
$ws.Range("M116").Value = 1186

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1697.8462
$ws.Range("I3").Value = 1108
$ws.Range("K3").Value = 1108
$ws.Range("M3").Value = -994
$ws.Range("H61").Value = 29300
$ws.Range("J61").Value = 29300
$ws.Range("L61").Value = 29300
$ws.Range("N61").Value = -29926

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 88.23529000000001
$ws.Range("I7").Value = 82.25
$ws.Range("J7").Value = 93.55556
$ws.Range("K7").Value = 82.25
$ws.Range("L7").Value = 93.55556
$ws.Range("M7").Value = 30.75
$ws.Range("N7").Value = -319.55556
$ws.Range("H26").Value = 30017.334
$ws.Range("J26").Value = 30017.334
$ws.Range("L26").Value = 30017.334
$ws.Range("N26").Value = -30591.334
$ws.Range("H31").Value = 7756426
$ws.Range("I31").Value = 5210.613
$ws.Range("J31").Value = 27780400
$ws.Range("K31").Value = 5210.613
$ws.Range("L31").Value = 27780400
$ws.Range("M31").Value = -4915.613
$ws.Range("N31").Value = -27780990
$ws.Range("H34").Value = 7756426
$ws.Range("I34").Value = 5210.613
$ws.Range("J34").Value = 27780400
$ws.Range("K34").Value = 5210.613
$ws.Range("L34").Value = 27780400
$ws.Range("M34").Value = -5008.613
$ws.Range("N34").Value = -27780804
$ws.Range("H132").Value = 27780186
$ws.Range("I132").Value = 50001936
$ws.Range("K132").Value = 150005808
$ws.Range("M132").Value = -150003278

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1375
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 1733.3334
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 5200.0002
$ws.Range("M25").Value = -731
$ws.Range("N25").Value = -5538.0002
$ws.Range("H30").Value = 1375
$ws.Range("I30").Value = 300
$ws.Range("J30").Value = 1733.3334
$ws.Range("K30").Value = 900
$ws.Range("L30").Value = 5200.0002
$ws.Range("M30").Value = -798
$ws.Range("N30").Value = -5404.0002
$ws.Range("H38").Value = 224.53334
$ws.Range("I38").Value = 791.5
$ws.Range("J38").Value = 137.3077
$ws.Range("K38").Value = 2374.5
$ws.Range("L38").Value = 411.9231
$ws.Range("M38").Value = -2027.5
$ws.Range("N38").Value = -1105.9231
$ws.Range("H107").Value = 1331.5294
$ws.Range("J107").Value = 1554.1428
$ws.Range("L107").Value = 4662.428400000001
$ws.Range("N107").Value = -8502.428400000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 19999.666
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 19999.666
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 19999.666
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -20553.666
$ws.Range("H107").Value = 1802.1177
$ws.Range("I107").Value = 2872.5557
$ws.Range("J107").Value = 597.875
$ws.Range("K107").Value = 2872.5557
$ws.Range("L107").Value = 597.875
$ws.Range("M107").Value = -952.5556999999999
$ws.Range("N107").Value = -4437.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1284.5264
$ws.Range("I22").Value = 575
$ws.Range("J22").Value = 1368
$ws.Range("K22").Value = 575
$ws.Range("L22").Value = 1368
$ws.Range("M22").Value = -280
$ws.Range("N22").Value = -1958
$ws.Range("H27").Value = 1284.5264
$ws.Range("I27").Value = 575
$ws.Range("J27").Value = 1368
$ws.Range("K27").Value = 575
$ws.Range("L27").Value = 1368
$ws.Range("M27").Value = -468
$ws.Range("N27").Value = -1582
$ws.Range("H46").Value = 736.9259
$ws.Range("I46").Value = 485.7143
$ws.Range("K46").Value = 485.7143
$ws.Range("M46").Value = -297.7143

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 903.2083
$ws.Range("I113").Value = 386.75
$ws.Range("J113").Value = 1936.125
$ws.Range("K113").Value = 1160.25
$ws.Range("L113").Value = 5808.375
$ws.Range("M113").Value = 1009.75
$ws.Range("N113").Value = -10148.375
